$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: F1 now holds the "image must be a link" note,
# and the old Multimedia_3 / Multimedia_4 columns (G1/H1) are removed.
$ws.Range("F1").Value = "(La imagen debe ser un link)"
$ws.Range("G1:H1").ClearContents()

# Widen column E slightly (target stored width ~17.57 chars; closest
# value the engine's pixel-quantized ColumnWidth can reproduce is 17.5)
$ws.Columns.Item(5).ColumnWidth = 16.71

# Update selection / view to match the saved state
$ws.Range("H9").Select()
